$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update raise target R (C2): 12 -> 25
$ws.Range("C2").Value = 25

# Update max Supply mS (C5): 1000000000 -> 5000000
$ws.Range("C5").Value = 5000000

# Update base Price floor f (C4): now a formula derived from the fixed
# floor price (1e9, scaled by contract decimals) divided by max supply
$ws.Range("C4").Formula = "=0.005*1000000000/C5"

# Add a note cell documenting the floor price used in the contract
$ws.Range("D29").Value = "floorprice=1000000000"

# restore selection like authors last save
$ws.Range("C3").Select()
